$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.275776982307434
$ws.Range("B1").Value = 2.19999623298645
$ws.Range("C1").Value = 4.727369785308838
$ws.Range("D1").Value = 3.183247327804565
$ws.Range("E1").Value = 1.360368728637695
